$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 29.576068633124976
$ws.Range("C2").Value = 16.626316213124994
$ws.Range("D2").Value = 24.048378633124969
$ws.Range("E2").Value = 24.913831213124979

# Row 3 data values
$ws.Range("B3").Value = 27.750552391249926
$ws.Range("C3").Value = 18.431477713124991
$ws.Range("D3").Value = 22.204002391249958
$ws.Range("E3").Value = 21.178336213124965

# Update selection to match the new selected range
$ws.Range("B1:E3").Select()
